$wb = $excel.ActiveWorkbook

# --- Config1: add Hawkeye / OPTGen rows for the "gromacs" benchmark ---
$ws1 = $wb.Worksheets.Item("Config1")

$ws1.Range("C29").Value = 50000001
$ws1.Range("D29").Value = 124667969
$ws1.Range("E29").Value = 30089
$ws1.Range("F29").Value = 15881
$ws1.Range("G29").Value = 14208

$ws1.Range("C30").Value = 50000001
$ws1.Range("D30").Value = 124667969
$ws1.Range("E30").Value = 441
$ws1.Range("F30").Value = 19
$ws1.Range("G30").Formula = "=E30-F30"
$ws1.Range("J30").Formula = "=F30/E30"

$ws1.Range("J31").Select() | Out-Null

# --- Config2: same new benchmark rows for the second cache configuration ---
$ws2 = $wb.Worksheets.Item("Config2")
$ws2.Activate()

$ws2.Range("C29").Value = 50000001
$ws2.Range("D29").Value = 124667224
$ws2.Range("E29").Value = 30134
$ws2.Range("F29").Value = 15919
$ws2.Range("G29").Value = 14215

$ws2.Range("C30").Value = 50000001
$ws2.Range("D30").Value = 124667224
$ws2.Range("E30").Value = 216
$ws2.Range("F30").Value = 18
$ws2.Range("G30").Formula = "=E30-F30"
$ws2.Range("J30").Formula = "=F30/E30"

$ws2.Range("I3").Select() | Out-Null
